$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12322
$ws1.Range("F7").Value = 12258
$ws1.Range("F15").Value = 3577

# Sheet "全部类型" (sheet4): same underlying rows, update matching values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12322
$ws4.Range("F8").Value = 12258
$ws4.Range("F17").Value = 3577
